# Weekly NYPD CompStat (88th Precinct) data refresh.
# Report period rolls forward one week (1/13-1/19/2025 -> 1/20-1/26/2025,
# "Number 3" -> "Number 4") and the Crime Complaints grid (rows 16-30) is
# refreshed with newly collected crime figures, recomputing the % Chg columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/issue number + reporting week dates ---
$ws.Range("A8").Value = "Volume 32   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/20/2025  Through  1/26/2025"

# --- Crime Complaints grid (rows 16-30): counts + computed % changes ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 28.571428571428
$ws.Range("I16").Value = 8
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = 60
$ws.Range("M16").Value = -61.904761904761
$ws.Range("N16").Value = -90.909090909090
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 133.333333333333
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 81.818181818181
$ws.Range("I17").Value = 19
$ws.Range("J17").Value = 9
$ws.Range("K17").Value = 111.111111111111
$ws.Range("L17").Value = 137.5
$ws.Range("M17").Value = 90
$ws.Range("N17").Value = -44.117647058823
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 14.285714285714
$ws.Range("I18").Value = 8
$ws.Range("J18").Value = 7
$ws.Range("K18").Value = 14.285714285714
$ws.Range("L18").Value = -46.666666666666
$ws.Range("M18").Value = -20
$ws.Range("N18").Value = -85.964912280701
$ws.Range("C19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 17
$ws.Range("G19").Value = 12
$ws.Range("H19").Value = 41.666666666666
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 12
$ws.Range("K19").Value = 16.666666666666
$ws.Range("L19").Value = -17.647058823529
$ws.Range("M19").Value = -17.647058823529
$ws.Range("N19").Value = -57.575757575757
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -80
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -58.823529411764
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -58.823529411764
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 75
$ws.Range("N20").Value = -88.333333333333
$ws.Range("C21").Value = 15
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 0
$ws.Range("G21").Value = 54
$ws.Range("H21").Value = 18.518518518518
$ws.Range("I21").Value = 59
$ws.Range("J21").Value = 50
$ws.Range("K21").Value = 18
$ws.Range("L21").Value = 7.272727272727
$ws.Range("M21").Value = -4.838709677419
$ws.Range("N21").Value = -78.700361010830
$ws.Range("C23").Value = 4
$ws.Range("D23").Value = 2
$ws.Range("D23").NumberFormat = "#,##0"
$ws.Range("E23").Value = 100
$ws.Range("E23").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F23").Value = 11
$ws.Range("H23").Value = 120
$ws.Range("I23").Value = 11
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = 175
$ws.Range("L23").Value = 266.666666666667
$ws.Range("M23").Value = -8.333333333333
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 42
$ws.Range("G24").Value = 35
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 38
$ws.Range("J24").Value = 34
$ws.Range("K24").Value = 11.764705882352
$ws.Range("L24").Value = -35.593220338983
$ws.Range("M24").Value = -40.625
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 200
$ws.Range("F25").Value = 13
$ws.Range("G25").Value = 13
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 12
$ws.Range("J25").Value = 11
$ws.Range("K25").Value = 9.090909090909
$ws.Range("L25").Value = -20
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 21
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = -22.222222222222
$ws.Range("L26").Value = -8.695652173913
$ws.Range("M26").Value = 40
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 3
$ws.Range("L28").Value = -25
$ws.Range("N29").Value = -77.777777777777
$ws.Range("N30").Value = -77.777777777777
